# Hoàn thiện Ngoại Trú
# Update the reception record's generated code/number (2049 -> 3013) and the
# derived IdCardNo / InsCardNo values on both the "Data" and "Check" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Data" ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("A2").Value = 3013
$wsData.Range("E2").Value = 46200608013
$wsData.Range("X2").Value = "DN4127460130013"

# --- Sheet "Check" ---
$wsCheck = $wb.Worksheets.Item("Check")
$wsCheck.Range("A2").Value = 3013
$wsCheck.Range("C2").Value = "DN4127460130013"
